$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Remember the sheet that was active before this script runs (Action1), so we
# can restore it at the end -- editing/selecting cells on "Global" must not
# change which tab is active in the saved workbook.
$origActive = $wb.ActiveSheet

# --- Header row (row 1): USERNAME / PASSWORD -------------------------------
$ws.Range("A1").Value = "USERNAME"
$ws.Range("B1").Value = "PASSWORD"

# --- Data row (row 2): sample credentials -----------------------------------
$ws.Range("A2").Value = 88996
$ws.Range("B2").Value = "Qatar@2021"

# --- Borders: thin box framing the A2:B2 input cells ------------------------
# Give A2 its top+bottom thin edges first...
$ws.Range("A2").Borders.Item(8).LineStyle = 1    # xlEdgeTop
$ws.Range("A2").Borders.Item(8).Weight = 2       # xlThin
$ws.Range("A2").Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$ws.Range("A2").Borders.Item(9).Weight = 2       # xlThin

# ...copy that border formatting onto B2 (so the two cells share the same
# top/bottom edge definition instead of building it twice)...
[void]$ws.Range("A2").Copy()
[void]$ws.Range("B2").PasteSpecial(-4122)        # xlPasteFormats

# ...then extend B2 with its own right edge so the pair reads as one boxed,
# two-cell input control.
$ws.Range("B2").Borders.Item(10).LineStyle = 1   # xlEdgeRight
$ws.Range("B2").Borders.Item(10).Weight = 2      # xlThin

# --- Column widths (best-fit to the new header/value text) ------------------
$ws.Columns.Item(1).ColumnWidth = 10.64
$ws.Columns.Item(2).ColumnWidth = 11.5

# --- Selection / view --------------------------------------------------------
[void]$ws.Range("B2").Select()

# Restore the originally active tab (Action1) so only the Global sheet's
# content/formatting changed, not which sheet is active.
[void]$origActive.Activate()
